$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.943.23"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +4.65%  '
$ws.Range("D3").Value = "'3.146.83"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +3.15%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'590.50"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.86%  '
$ws.Range("D6").Value = "'147.54"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +3.76%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = "'3.140.76"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +3.35%  '
$ws.Range("D9").Value = "'0.539"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.84%  '
$ws.Range("D10").Value = "'0.165"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +20.65%  '
$ws.Range("D11").Value = "'5.74"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +5.92%  '
$ws.Range("E12").Value = '  +2.14%  '
$ws.Range("D13").Value = "'0.0000258"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +9.58%  '
$ws.Range("D14").Value = "'35.97"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +3.84%  '
$ws.Range("E15").Value = '  +0.40%  '
$ws.Range("D16").Value = "'3.668.84"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +3.22%  '
$ws.Range("D17").Value = "'63.874.96"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +4.51%  '
$ws.Range("D18").Value = "'3.148.29"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +3.31%  '
$ws.Range("D19").Value = "'7.16"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.44%  '
$ws.Range("D20").Value = "'470.56"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +5.77%  '
$ws.Range("D21").Value = "'14.26"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").Value = "'0.734"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.63%  '
$ws.Range("D23").Value = "'7.55"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +4.44%  '
$ws.Range("D24").Value = "'13.38"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.50%  '
$ws.Range("D25").Value = "'82.84"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.38%  '
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("D27").Value = "'8.64"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +7.32%  '
$ws.Range("D28").Value = "'2.71"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +3.71%  '
$ws.Range("E29").Value = '  +0.21%  '
$ws.Range("E30").Value = '  -2.63%  '
$ws.Range("D31").Value = "'6.86"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +6.82%  '
$ws.Range("D32").Value = "'27.13"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.82%  '
$ws.Range("E33").Value = '  +2.63%  '
$ws.Range("D34").Value = "'0.0₃0875"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +8.85%  '
$ws.Range("D35").Value = "'2.42"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +11.82%  '
$ws.Range("E36").Value = '  +2.97%  '
$ws.Range("D37").Value = "'3.40"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +15.86%  '
$ws.Range("D38").Value = "'6.16"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.15%  '
$ws.Range("D39").Value = "'50.96"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.84%  '
$ws.Range("D40").Value = "'450.02"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +9.70%  '
$ws.Range("D41").Value = "'8.72"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.51%  '
$ws.Range("D42").Value = "'0.0376"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +4.15%  '
$ws.Range("D43").Value = "'2.920.53"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +4.99%  '
$ws.Range("E44").Value = '  +7.56%  '
$ws.Range("D45").Value = "'0.112"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +4.95%  '
$ws.Range("D46").Value = "'2.16"
$ws.Range("D46").ClearFormats()
$ws.Range("D47").Value = "'125.74"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.17%  '
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("E49").Value = '  +1.38%  '
$ws.Range("D50").Value = "'24.90"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +4.03%  '
$ws.Range("D51").Value = "'33.93"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -9.24%  '
